# Simulated Wild Card round and logged it.
# Adds the newly-signed/rostered running back "M.Sargent" as a new row
# at the bottom of the RB player-stats sheet, with all stat columns
# initialized to 0 (same shape as the other rows on this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

# New player row goes right after the last existing row (row 7 -> row 8).
$newRow = 8

$ws.Cells.Item($newRow, 1).Value = "M.Sargent"

# Stat columns B:J (PCOMPs/Rec, RecYDs, RecTDs, FMBLs, RATTs, RYDs, RTDs, 2PCs, FPTS)
# all start at 0, just like every other player row on this sheet.
$ws.Range("B" + $newRow + ":J" + $newRow).Value = 0

# Move/update the active selection to the next empty row below the new
# entry, mirroring Excel's normal "select next cell after data entry"
# behavior.
$ws.Range("J" + ($newRow + 1)).Select()
